# Apply forecast summary updates per diff: "Optuna Attempt (go back with original)"

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates (Seasonality Index column L, Inventory Coverage column H, MyForecast column D) ---

$wsForecast.Range("L2").Value = 0.98
$wsForecast.Range("L3").Value = 1.04
$wsForecast.Range("L4").Value = 0.82
$wsForecast.Range("L5").Value = 1.01

$wsForecast.Range("H6").Value = 16.56

$wsForecast.Range("H7").Value = 16.97
$wsForecast.Range("L7").Value = 0.84

$wsForecast.Range("H8").Value = 15.97
$wsForecast.Range("L8").Value = 1.04

$wsForecast.Range("H9").Value = 14.97
$wsForecast.Range("L9").Value = 0.97

$wsForecast.Range("H10").Value = 18.28
$wsForecast.Range("L10").Value = 0.88

$wsForecast.Range("D11").Value = 1
$wsForecast.Range("H11").Value = 21.36
$wsForecast.Range("L11").Value = 1.13

$wsForecast.Range("D12").Value = 1
$wsForecast.Range("H12").Value = 18.87
$wsForecast.Range("L12").Value = 0.89

$wsForecast.Range("D13").Value = 1
$wsForecast.Range("H13").Value = 17.87
$wsForecast.Range("L13").Value = 0.87

$wsForecast.Range("D14").Value = 1
$wsForecast.Range("H14").Value = 18.2
$wsForecast.Range("L14").Value = 1.04

$wsForecast.Range("D15").Value = 1
$wsForecast.Range("H15").Value = 17.2
$wsForecast.Range("L15").Value = 0.88

$wsForecast.Range("D16").Value = 1
$wsForecast.Range("H16").Value = 16.2
$wsForecast.Range("L16").Value = 0.96

$wsForecast.Range("D17").Value = 1
$wsForecast.Range("H17").Value = 15.2

# --- Summary sheet updates (text values) ---
# These cells are stored as text strings (not numbers) in the original file,
# so force text formatting before assigning to avoid Excel auto-converting
# the numeric-looking strings into numeric cells.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "38"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "23"
